# Generate Report for Archive
#
# 1. Change the status text "Ready for handoff" -> "In Translation" everywhere
#    it appears (Overview sheet columns E/F rows 2-3, and the zh-cn / de-de
#    sheets column C rows 2-3).
# 2. Shrink the corresponding "Status" columns that auto-fit to that text,
#    since the new text is shorter than the old one.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Columns E and F were sized to fit "Ready for handoff"; re-fit to the
# shorter replacement text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
